$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 80

# Columns A-D (Date, Time, Weekday, Week) hold plain text in this sheet
# (e.g. "2023-06-26", "22:28:02", "26"), not real dates/numbers. Force a
# text number format before assigning so Excel's automatic value
# inference doesn't turn them into a date serial / number, then clear the
# temporary formatting so the new row keeps the sheet's default styling.
$textRange = $ws.Range("A" + $row + ":D" + $row)
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2023-06-26"
$ws.Cells.Item($row, 2).Value = "22:28:02"
$ws.Cells.Item($row, 3).Value = "Monday"
$ws.Cells.Item($row, 4).Value = "26"

$textRange.ClearFormats()

# Columns E-T (city resale numbers) are plain numeric values.
$numCols = @{
    5  = 122810
    6  = 134247
    7  = 163695
    8  = 133563
    9  = 177282
    10 = 115113
    11 = 203598
    12 = 226276
    13 = 176277
    14 = 104330
    15 = 39685
    16 = 33786
    17 = 52196
    18 = -1
    19 = 36049
    20 = -1
}

foreach ($col in $numCols.Keys) {
    $ws.Cells.Item($row, $col).Value = $numCols[$col]
}
